$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '65.848.76'
$ws.Cells.Item(2, 5).Value = '  +6.90%  '

$ws.Cells.Item(3, 4).Value = '3.015.13'
$ws.Cells.Item(3, 5).Value = '  +4.10%  '

$ws.Cells.Item(4, 5).Value = '  +0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '585.63'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +3.10%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '154.81'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +7.81%  '

$ws.Cells.Item(7, 5).Value = '  -0.07%  '

$ws.Cells.Item(8, 4).Value = '3.010.90'
$ws.Cells.Item(8, 5).Value = '  +4.03%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.518'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +2.52%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.00'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.34%  '

$ws.Cells.Item(11, 5).Value = '  +4.82%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.451'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +4.67%  '

$ws.Cells.Item(13, 5).Value = '  +3.28%  '

$ws.Cells.Item(14, 5).Value = '  +7.70%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.126'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.76%  '

$ws.Cells.Item(16, 4).Value = '65.826.46'
$ws.Cells.Item(16, 5).Value = '  +6.86%  '

$ws.Cells.Item(17, 4).Value = '3.517.03'
$ws.Cells.Item(17, 5).Value = '  +4.21%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.97'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +6.11%  '

$ws.Cells.Item(19, 4).Value = '3.018.60'
$ws.Cells.Item(19, 5).Value = '  +4.30%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '460.88'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +6.38%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.78'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +5.37%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.685'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +4.11%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '7.37'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +7.89%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '82.05'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +3.44%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '12.61'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +5.68%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.25'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +11.72%  '

$ws.Cells.Item(27, 5).Value = '  +8.36%  '

$ws.Cells.Item(28, 5).Value = '  -0.07%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.43'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +19.07%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '7.95'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +13.17%  '

$ws.Cells.Item(31, 2).Value = 'PEPE'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.0000104'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -3.00%  '

$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.61'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +4.21%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.112'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +5.48%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '27.00'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +5.72%  '

$ws.Cells.Item(35, 5).Value = '  -0.13%  '

$ws.Cells.Item(36, 5).Value = '  +3.50%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.80'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +7.50%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.16'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +11.62%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.01'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +6.51%  '

$ws.Cells.Item(40, 2).Value = 'OKB'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '49.46'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.29%  '

$ws.Cells.Item(41, 2).Value = 'Arweave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '45.06'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +13.67%  '

$ws.Cells.Item(42, 2).Value = 'TheGraph'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.303'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +13.94%  '

$ws.Cells.Item(43, 2).Value = 'Kaspa'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.121'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +6.28%  '

$ws.Cells.Item(44, 5).Value = '  +3.10%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '389.36'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +11.63%  '

$ws.Cells.Item(46, 4).Value = '2.799.15'
$ws.Cells.Item(46, 5).Value = '  +3.96%  '

$ws.Cells.Item(47, 5).Value = '  +5.21%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '134.94'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +1.54%  '

$ws.Cells.Item(49, 5).Value = '  -0.05%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '23.69'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +9.58%  '
